$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column C for rows 2-7 from 45212 to 45221
for ($row = 2; $row -le 7; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45212) {
        $cell.Value = 45221
    }
}
